# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - Home (row 2) target depth stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 349
$wsOff.Range("C2").Value = 246
$wsOff.Range("D2").Value = 81
$wsOff.Range("E2").Value = 38

# DEF sheet - Home (row 2) target depth stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 415
$wsDef.Range("C2").Value = 287
$wsDef.Range("D2").Value = 108
$wsDef.Range("E2").Value = 49
$wsDef.Range("F2").Value = 3
